# direcciones add y update
$wb = $excel.ActiveWorkbook

# --- Sheet "Deuda": move selection from B8 to A8 ---
$wsDeuda = $wb.Worksheets.Item("Deuda")
$wsDeuda.Range("A8").Select()

# --- Sheet "horas presenciales": add new rows 10-14 and move selection ---
$wsHoras = $wb.Worksheets.Item("horas presenciales")

# Seed the new date cells (A10:A14) with the same formatting as the
# existing date column by copying an already-formatted cell first, then
# overwrite just the values.
$wsHoras.Range("A8").Copy($wsHoras.Range("A10"))
$wsHoras.Range("A8").Copy($wsHoras.Range("A11"))
$wsHoras.Range("A8").Copy($wsHoras.Range("A12"))
$wsHoras.Range("A8").Copy($wsHoras.Range("A13"))
$wsHoras.Range("A8").Copy($wsHoras.Range("A14"))

$wsHoras.Range("A10").Value = 42088
$wsHoras.Range("B10").Value = 6
$wsHoras.Range("C10").Value = 6
$wsHoras.Range("D10").Value = 6
$wsHoras.Range("E10").Value = 0

$wsHoras.Range("A11").Value = 42089
$wsHoras.Range("B11").Value = 3
$wsHoras.Range("C11").Value = 0
$wsHoras.Range("D11").Value = 3
$wsHoras.Range("E11").Value = 0

$wsHoras.Range("A12").Value = 42090
$wsHoras.Range("B12").Value = 6
$wsHoras.Range("C12").Value = 6
$wsHoras.Range("D12").Value = 0
$wsHoras.Range("E12").Value = 0

$wsHoras.Range("A13").Value = 42100
$wsHoras.Range("B13").Value = 3
$wsHoras.Range("C13").Value = 0
$wsHoras.Range("D13").Value = 0
$wsHoras.Range("E13").Value = 3

$wsHoras.Range("A14").Value = 42100
$wsHoras.Range("B14").Value = 4
$wsHoras.Range("C14").Value = 0
$wsHoras.Range("D14").Value = 0
$wsHoras.Range("E14").Value = 4

$wsHoras.Range("B16").Select()
